$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows below had their match/odds data swapped between each pair while
# keeping the row index (column A) and the Date (column D) fixed in place.
# Swap everything from column B through column AD for each pair.
$pairs = @(
    @(192, 193),
    @(256, 257),
    @(270, 271)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $rangeA = $ws.Range("B$r1`:AD$r1")
    $rangeB = $ws.Range("B$r2`:AD$r2")

    $valsA = $rangeA.Value2
    $valsB = $rangeB.Value2

    $rangeA.Value2 = $valsB
    $rangeB.Value2 = $valsA
}
